# sedani/chassis_control_board/control_board_3.1_parts.xlsx
# "added flyback diode and worked on routing"
#
# T1's relay driver part changes from an NPN BJT to an N-channel MOSFET,
# and a new D4 flyback diode row is added to the parts table, along with
# its Digikey part number.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (T1): swap the NPN transistor for an N MOSFET and update its part number.
# Write the B2/C2 pair first so new shared strings are interned in commit order.
$ws.Range("B2").Value = "N MOSFET for relay"
$ws.Range("C2").Value = "NTR4003NT1GOSCT-ND"

# Row 4 (new): D4, the flyback diode added for the relay.
$ws.Range("B4").Value = "flyback diode for relay"
$ws.Range("C4").Value = "1655-1502-1-ND"
$ws.Range("A4").Value = "D4"

# Row 3: the voltage regulator row now also gets its Digikey part number filled in.
$ws.Range("C3").Value = "AZ2940D-5.0TRE1DICT-ND"

# Column C widened slightly to fit the longer part numbers.
$ws.Columns("C").ColumnWidth = 22.33

# Leave the active selection on C3, matching where editing finished.
$ws.Range("C3").Select()
